# "updated UI for batch operation"
# - rename shared string used by the Task ID column (B2:B3) from
#   "acma_check" to "inter_transmission_merge"
# - move the active cell selection on Sheet1 from D6 to B8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the task name text (shared string) used by both data rows.
$ws.Range("B2").Value = "inter_transmission_merge"
$ws.Range("B3").Value = "inter_transmission_merge"

# Update the UI selection state to the new active cell.
$ws.Range("B8").Select() | Out-Null
